$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder header row: B1/C1/D1 become Phone Number / Nationality / Email
$ws.Range("B1").Value = "Phone Number"
$ws.Range("C1").Value = "Nationality"
$ws.Range("D1").Value = "Email"

# Move the email values from column B to column D for rows 2 and 3,
# clearing out the now-empty B/C cells along the way.
$ws.Range("D2").Value = $ws.Range("B2").Value2
$ws.Range("B2").ClearContents()

$ws.Range("D3").Value = $ws.Range("B3").Value2
$ws.Range("B3").ClearContents()
